$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E1").Value = 3.5396999999999998
$ws.Range("F1").Value = 98.522300000000001

$ws.Range("E2").Value = 3.6261999999999999
$ws.Range("F2").Value = 98.560299999999998

$ws.Range("E3").Value = 5.7118000000000002
$ws.Range("F3").Value = 98.6905

$ws.Range("F3").Select()
